$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "ReceivesAction" row (row 12) entirely
$ws.Rows.Item(12).Delete()

# Reword a few templates in column B
$ws.Range("B8").Value = "It has <HasA>"
$ws.Range("B5").Value = "It likes <Desires>"
$ws.Range("B4").Value = "It hates <NotDesires>"

# Replace the "Example cat" column (C) with an "Expected POS" column
$ws.Range("C1").Value = "Expected POS"
$ws.Range("C2").Value = "NOUN"
$ws.Range("C3").Value = "VERB"
$ws.Range("C4").Value = "NOUN, VERB"
$ws.Range("C5").Value = "NOUN, VERB"
$ws.Range("C6").Value = "VERB"
$ws.Range("C7").Value = "NOUN"
$ws.Range("C8").Value = "NOUN"
$ws.Range("C9").Value = "NOUN, VERB, ADJ"
$ws.Range("C10").Value = "NOUN"
$ws.Range("C11").Value = "NOUN"

# Match the saved selection state of the source workbook
$ws.Range("C11").Select()
